# Generate Report for Handback
#
# Populates the handback-report row for file
# "5caf57f1-57b5-42a5-822c-8f5fac5aa3c1" on both locale sheets (zh-cn, de-de):
#   - Latest Target File (I7)      -> the source .md hyperlink (new)
#   - Latest Handback File (J7)    -> the locale-specific xlf that was handed back
#   - Latest Handback DateTime(K7) -> the timestamp the handback was processed
#   - Error Detail (P7)            -> handback version-mismatch message
#
# The handback file used to generate this report was not built from the
# latest handoff, so the status message records both commit SHAs.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4a422da47b425b06ec21e2ef529392bb05896f13/e2e/5caf57f1-57b5-42a5-822c-8f5fac5aa3c1.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/236c5b64e3c88027d9feee8c0165e17e09edc83f/e2e/5caf57f1-57b5-42a5-822c-8f5fac5aa3c1.md."

$mdDisplay = "5caf57f1-57b5-42a5-822c-8f5fac5aa3c1.md"

# zh-cn sheet
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("J7").Value = "5caf57f1-57b5-42a5-822c-8f5fac5aa3c1.8be61c10765e08b6ee48da4015b3bd335d50faab.zh-cn.xlf"
$wsZh.Range("K7").Value = "2016-08-15 12:50:53"
$wsZh.Range("P7").Value = $errorDetail
$wsZh.Hyperlinks.Add($wsZh.Range("I7"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/236c5b64e3c88027d9feee8c0165e17e09edc83f/e2e/5caf57f1-57b5-42a5-822c-8f5fac5aa3c1.md", [type]::Missing, [type]::Missing, $mdDisplay)

# de-de sheet
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("J7").Value = "5caf57f1-57b5-42a5-822c-8f5fac5aa3c1.8be61c10765e08b6ee48da4015b3bd335d50faab.de-de.xlf"
$wsDe.Range("K7").Value = "2016-08-15 12:51:01"
$wsDe.Range("P7").Value = $errorDetail
$wsDe.Hyperlinks.Add($wsDe.Range("I7"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/236c5b64e3c88027d9feee8c0165e17e09edc83f/e2e/5caf57f1-57b5-42a5-822c-8f5fac5aa3c1.md", [type]::Missing, [type]::Missing, $mdDisplay)
